$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bValues = @(0.71875, 0.71875, 0.6875, 0.625, 0.59375, 0.5625, 0.53125, 0.515625, 0.46875, 0.421875, 0.40625, 0.453125, 0.4375, 0.375, 0.375, 0.421875, 0.359375, 0.546875, 0.390625, 0.375, 0.359375, 0.375, 0.359375, 0.375, 0.359375, 0.375, 0.359375, 0.40625, 0.390625, 0.390625, 0.390625, 0.390625, 0.390625, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.359375, 0.359375, 0.359375, 0.359375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.375, 0.46875, 0.46875, 0.640625, 0.46875, 0.546875, 0.65625, 0.46875, 0.515625, 0.5, 0.546875, 0.484375, 0.5625, 0.5)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

$null = $ws.Range("A2:B115").Select()
